$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New, re-sorted and extended data table (rows 2-12)
$data = @(
    @("Basic CNN",     "No",  200, "Rotation", 0.61419999999999997),
    @("Basic CNN",     "No",  300, "Rotation", 0.72199999999999998),
    @("Basic CNN",     "No",  300, "All",      0.77410000000000001),
    @("Advanced CNN",  "No",  200, "Rotation", 0.6804),
    @("Advanced CNN",  "Yes", 200, "Rotation", 0.62419999999999998),
    @("Advanced CNN",  "Yes", 300, "Rotation", 0.81330000000000002),
    @("Advanced CNN",  "Yes", 300, "All",      0.76880000000000004),
    @("Advanced CNN",  "No",  300, "Rotation", 0.8216),
    @("Advanced CNN",  "No",  400, "Rotation", 0.80110000000000003),
    @("Advanced CNN",  "No",  300, "All",      0.84150000000000003),
    @("Efficient Net", "No",  300, "All",      0.98650000000000004)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}

$ws.Range("L9").Select()
